$wb = $excel.ActiveWorkbook

# Add the two new sheets after o_10, in order
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "o_20"
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $ws2)
$ws3.Name = "o_20_jumbled"

$promptA = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0 0
 G 0 0 1 0 0 1 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 1 0 0 0 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 1 0 0 1 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0

Solution: A -> E -> F -> G -> H -> L -> P
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node P? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 1 0 0 0 0 1 0 0 1 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 0 0 0 1 0 0 0 0
 I 0 0 0 0 1 0 0 0 0 1 0 0 1 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 K 0 0 0 0 0 0 1 0 0 1 0 1 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 1
 M 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
    
"@
$promptB = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 24 nodes labelled A to X. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0

Solution: A -> F -> G -> L -> M -> N -> O -> T -> Y
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node X? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 1 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@
$promptC = @"
 Given is the adjacency matrix for a unweighted undirected graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0

Solution: A -> F -> G -> L -> M -> N -> O -> T -> Y
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node Y? Return the sequence of nodes in response.

   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 1 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 1 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0 0 0 1 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 1
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 1
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0
    
"@

# --- Sheet1 (o_10): add column E header + update row2 ---
$ws1.Range("D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("A2").Value = $promptA
$ws1.Range("B2").Value = "A -> E -> I -> J -> K -> L -> P"
$ws1.Range("C2").Value = "The shortest path from node A to node P is: A -> B -> F -> E -> I -> J -> K -> L -> P."
$ws1.Range("D2").Value = "invalid input"
$ws1.Range("E2").Value = "1/7"

# --- Sheet2 (o_20): header row + row2 values ---
$ws1.Range("A1:E1").Copy($ws2.Range("A1:E1"))
$ws2.Range("A2").Value = $promptB
$ws2.Range("B2").Value = "A -> E -> F -> G -> H -> I -> N -> S -> X"
$ws2.Range("C2").Value = "The shortest path from node A to node X is: A -> E -> J -> O -> S -> X."
$ws2.Range("D2").Value = "invalid input"
$ws2.Range("E2").Value = "5/9"

# --- Sheet3 (o_20_jumbled): header row + row2 values ---
$ws1.Range("A1:E1").Copy($ws3.Range("A1:E1"))
$ws3.Range("A2").Value = $promptC
$ws3.Range("B2").Value = "A -> F -> G -> H -> M -> N -> O -> T -> Y"
$ws3.Range("C2").Value = "The shortest path from node A to node Y is: A -> F -> G -> H -> M -> N -> O -> T -> Y"
$ws3.Range("D2").Value = "invalid input"
$ws3.Range("E2").Value = "9/9"

# Restore original active sheet/tab selection
$ws1.Select()
